$d = $word.ActiveDocument

# 1) Remove the stray _GoBack bookmark that sits alone in the empty
#    paragraph right before the "Buttons" heading.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Fix the duplicated phrase in the Mealy-machine sentence.
$d.Content.Find.Execute(
    "can help visualize the Mealy machine can help visualize the state transitions:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "can help visualize the Mealy machine state transitions:", 2) | Out-Null

# 3) Re-insert the _GoBack bookmark at the point where the duplicate text
#    used to be edited, splitting the run exactly like Word does when the
#    last edit position is remembered.
$r = $d.Content
$r.Find.Execute(
    "A finite state diagram can help visualize the Mealy machine ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$insPoint = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $insPoint)
